# Notes regarding programming MCU
# Adds a "." run to the final paragraph (closing the sentence) and then
# appends a set of new paragraphs documenting "Revision 1" and "Revision 2"
# of the ICSP / USB-SPI programming notes.

$d = $word.ActiveDocument

# --- 1. Close out the existing last paragraph with a "." ---------------
$lastPara = $d.Paragraphs.Last
$r = $lastPara.Range
$r.Collapse(0)
$r.InsertAfter(".")

# --- helper-ish pattern: always grab the doc's current last paragraph,
#     insert a new paragraph after it, then fill in its text/style -------

# Heading3: "Revision 1:"
$d.Paragraphs.Last.Range.InsertParagraphAfter()
$p = $d.Paragraphs.Last
$p.Style = $d.Styles("Heading 3")
$p.Range.Text = "Revision 1:"

# Body paragraph
$d.Paragraphs.Last.Range.InsertParagraphAfter()
$p = $d.Paragraphs.Last
$p.Style = $d.Styles("Normal")
$p.Range.Text = "Went full circle essentially. Probably going to borrow a design off the internet to use the programmer. If the IC states that it supports USB 2.0, it works with all future versions as well (USB 3, 4). "

# Body paragraph
$d.Paragraphs.Last.Range.InsertParagraphAfter()
$p = $d.Paragraphs.Last
$p.Style = $d.Styles("Normal")
$p.Range.Text = "There are three protocols to programming a chip UART programming, SPI, and UDPI"

# Body paragraph
$d.Paragraphs.Last.Range.InsertParagraphAfter()
$p = $d.Paragraphs.Last
$p.Style = $d.Styles("Normal")
$p.Range.Text = "UART uses 2, SPI uses 1 ish and UDPI uses 1"

# Body paragraph
$d.Paragraphs.Last.Range.InsertParagraphAfter()
$p = $d.Paragraphs.Last
$p.Style = $d.Styles("Normal")
$p.Range.Text = "UDPI is the newest version of the programming protocols and I have no idea how it works"

# Heading3: "Revision 2:"
$d.Paragraphs.Last.Range.InsertParagraphAfter()
$p = $d.Paragraphs.Last
$p.Style = $d.Styles("Heading 3")
$p.Range.Text = "Revision 2:"

# Body paragraph
$d.Paragraphs.Last.Range.InsertParagraphAfter()
$p = $d.Paragraphs.Last
$p.Style = $d.Styles("Normal")
$p.Range.Text = "IVE BEEN LOSING MY MIND OVER THIS FOR 3 HOURS AND I THINK I GOT IT"

# Body paragraph
$d.Paragraphs.Last.Range.InsertParagraphAfter()
$p = $d.Paragraphs.Last
$p.Style = $d.Styles("Normal")
$p.Range.Text = "Ecosystem:"

# Body paragraph (the long "Convert data..." one, followed by the
# ". Wtf this took so long :sob:" tail appended onto the same paragraph)
$d.Paragraphs.Last.Range.InsertParagraphAfter()
$p = $d.Paragraphs.Last
$p.Style = $d.Styles("Normal")
$p.Range.Text = "Convert data from USB to SPI using an interface chip. Feed this data into a separate MCU, which contains the bootloader for the actual MCU. (Bootloader being the software that runs every time the power gets reset). This bootloader holds the software that facilitates the communication between the two devices using a protocol that is understood by the main MCU and bam coding done"
$tailRange = $p.Range
$tailRange.Collapse(0)
$tailRange.InsertAfter(". Wtf this took so long :sob:")
